$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix E9 / E10 -----------------------------------------------------
# They were stored as inline text "543237" / "526371"; the breakout run
# now writes the bsecode as a genuine number for these two existing rows.
$ws.Range("E9").Value = 543237
$ws.Range("E10").Value = 526371

# --- Append rows 11-13 (new screener snapshot, 10/06/2024 09:46:26) ---
$ws.Range("A11").Value = "10/06/2024 09:46:26"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "MAZDOCK"
$ws.Range("D11").Value = "Mazagon Dock Shipbuilders Ltd"
$ws.Range("F11").Value = -0.98
$ws.Range("G11").Value = 3122.1
$ws.Range("H11").Value = 1255452

$ws.Range("A12").Value = "10/06/2024 09:46:26"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "NMDC"
$ws.Range("D12").Value = "Nmdc Limited"
$ws.Range("F12").Value = -1.41
$ws.Range("G12").Value = 254.85
$ws.Range("H12").Value = 4873940

$ws.Range("A13").Value = "10/06/2024 09:46:26"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "GAIL"
$ws.Range("D13").Value = "Gail (india) Limited"
$ws.Range("F13").Value = -1.37
$ws.Range("G13").Value = 209.79
$ws.Range("H13").Value = 32925648

# bsecode (column E) on the new rows keeps the *text* representation seen
# on the older rows (e.g. "543237" rather than the number 543237), so
# write it through TEXT()+PasteSpecial(values) instead of .Value = "..."
# (a plain string assignment of a numeric-looking literal gets silently
# re-typed to a number by Excel, same as typing it into the grid).
$ws.Range("Z1").Formula = "=TEXT(543237,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("Z1").Formula = "=TEXT(526371,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("Z1").Formula = "=TEXT(532155,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
